$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"

$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 30.03.2025"

$ws.Range("B6").Value = "01.04."
$ws.Range("C6").Value = "02.04."
$ws.Range("D6").Value = "RECHNUNG VODAFONE GMBH 26432426"
$ws.Range("E6").Value = "41,37-"

$ws.Range("B7").Value = "05.04."
$ws.Range("C7").Value = "06.04."
$ws.Range("D7").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 50704860"
$ws.Range("E7").Value = "84,13-"

$ws.Range("B8").Value = "06.04."
$ws.Range("C8").Value = "07.04."
$ws.Range("D8").Value = "BEITRAG Allianz SE K-3256986"
$ws.Range("E8").Value = "57,52-"

$ws.Range("B9").Value = "09.04."
$ws.Range("C9").Value = "10.04."
$ws.Range("D9").Value = "PAYPAL JKTGJV"
$ws.Range("E9").Value = "71,98-"

$ws.Range("B10").Value = "11.04."
$ws.Range("C10").Value = "12.04."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,02-"

$ws.Range("D12").Value = "KONTOSTAND AM 16.04.2025"
$ws.Range("E12").Value = "280,02-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 22.04.2025"
